$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a value that must be stored as TEXT even if it looks numeric,
# by prefixing with an apostrophe (standard Excel "force text" input trick).
function Set-TextValue($rangeAddr, $text) {
    if ($text -eq "") {
        $ws.Range($rangeAddr).Value = ""
    } else {
        $ws.Range($rangeAddr).Value = "'" + $text
    }
}

# ---------- Row 8 ----------
$ws.Range("C8").Value = 56
Set-TextValue "D8" "1.0"
$ws.Range("E8").Value = "Rewiring of light point/ fan point/ exhaust fan point/ call bell point with 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper conductor 1.1 kV grade  and 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper earth conductor 1.1 kV grade (IS:694) in recessed  ISI marked MMS ( IS:9537 P - III ) virgin material PVC conduit & it's ISI marked (IS:3419-1988) accessories, round tiles, 1.2 mm thick MS box with earth terminal, 6 A switch, 3 pin ceiling rose/holder / 3 way connector , 3.0 mm thick ISI marked (IS:2036-1995) phenolic laminated sheet, Al.alloy / Cadmium plated iron/  brass  screws, cup washers, making connections, testing etc. as required. For specification of copper  Conductor,  Phenolic Laminated sheet's & Electrical/ Wiring accessories refer Chapter E - 04, E - 05 & E - 07 For additional technical parameters of product / work refer Annexure 'A' attached with this BSR"

# ---------- Row 9 ----------
$ws.Range("C9").Value = 66
Set-TextValue "D9" "4"
$ws.Range("E9").Value = "Long point  (up to 10 mtr.)"
$ws.Range("F9").Value = 662
Set-TextValue "G9" "43692.00"

# ---------- Row 10 ----------
Set-TextValue "A10" ""
$ws.Range("C10").Value = 26
Set-TextValue "D10" "2.0"
$ws.Range("E10").Value = "Rewiring of 3/5 pin 6 amp. Light plug point with 1.5 sq. mm nominal size  FR PVC insulated unsheathed flexible copper conductor 1.1 kV grade  and 1.5 sq. mm nominal size FR PVC insulated unsheathed flexible copper earth conductor 1.1 kV grade(IS:694)   in recessed ISI marked MMS ( IS:9537 P - III ) virgin material  PVC conduit & it's  ISI marked (IS:3419-1988) accessories, 1.2 mm thick  MS box with earth terminal of required size,  6 A  switch, 3/5 pin 6 A socket, 3.0 mm thick ISI marked (IS:2036-1995) phenolic laminated sheet, Al.alloy / Cadmium plated iron/ brass  screws, cup washers, making connections, testing etc. as required.  For specification of copper  Conductor,  Phenolic Laminated sheet's & Electrical/ Wiring accessories refer Chapter E - 04, E - 05 & E - 07 For additional technical parameters of product / work refer Annexure 'A' attached with this BSR"
$ws.Range("F10").Value = 0
Set-TextValue "G10" "0.00"

# ---------- Row 11 ----------
$ws.Range("C11").Value = 41
Set-TextValue "D11" "4.0"
$ws.Range("E11").Value = "P & F ISI marked (IS :3854) 16 amp. flush type non modular switch CM/L no. printed and made out from industrial grade Polycarbonate or fire resistant ABS material including cutting hole in tile and making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure `"A`" attached with this BSR ."
$ws.Range("F11").Value = 50
Set-TextValue "G11" "2050.00"

# ---------- Row 12 ----------
$ws.Range("C12").Value = 51
Set-TextValue "D12" "5.0"
$ws.Range("E12").Value = "Providing & Fixing of  of 3/5 pin 6 amp. flush type  non modular socket  made out from industrial grade Polycarbonate or fire resistant ABS material, brass terminal with Porcelain based back cover & captive screws including cutting hole in tile and making connection testing etc. as required.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure `"A`" attached with this BSR ."
$ws.Range("F12").Value = 33
Set-TextValue "G12" "1683.00"

# ---------- Row 13 ----------
$ws.Range("C13").Value = 84
Set-TextValue "G13" "4704.00"

# ---------- Row 14 ----------
Set-TextValue "A14" "Set"
$ws.Range("C14").Value = 58
Set-TextValue "D14" "13.0"
$ws.Range("E14").Value = "Plate Earthing  as per IS:3043 with Hot dipped G.I. Earth plate of size 600mm x 600mm x 6.0mm by embodying  3 to 4 mtr. below the ground level with 20  mm dia. G.I. 'B' class watering Pipe ,including all accessories like nut, bolts, reducer, nipple, wire meshed funnel, and Heavy duty weather proof poly-propylene earth pit chamber with lockable Jam free lid suitable for safe working load 5000 Kg or more of size Top Dia. 225 to 260 mm, Bottom Dia 300 to 350 mm. and Height  250 to 300 mm. and embodying the pipe  complete with alternate layers salt and coke/ charcoal, testing of earth resistance for value of 5 ohms or less  as required & must record by engineer in charge during site visit and ensure to enter in measurment book.All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure `"A`" attached with this BSR .   "
$ws.Range("F14").Value = 5733
Set-TextValue "G14" "332514.00"

# ---------- Row 15 ----------
Set-TextValue "A15" "Mtr."
$ws.Range("C15").Value = 14
Set-TextValue "D15" "23"
$ws.Range("E15").Value = "8 SWG G.I. ( Hot Dipped  ) Wire "
$ws.Range("F15").Value = 20
Set-TextValue "G15" "280.00"

# ---------- Row 16 ----------
Set-TextValue "A16" ""
$ws.Range("C16").Value = 2
Set-TextValue "D16" "15.0"
$ws.Range("E16").Value = "Providing & Fixing of  BEE  Star rated copper wounded double ball bearing capacitor start, aluminium body & Metallic  blade ceiling  fan  Conforming to all the performance requirements laid down in IS 374:2019 including all amendments, as applicable ; & Carry BIS licensing (i.e. ISI marking) with down rod up to 80 cm with secondary support safety cable ( steel rope) , cotter pin with 3 x 1.5 sq.mm pvc insulated flexible copper conductor making connection testing etc. as required. All as per pre approved by Engineer in charge.  For additional technical parameters of products/ work  , refer   Annexure `"A`" attached with this BSR ."
$ws.Range("F16").Value = 0
Set-TextValue "G16" "0.00"

# ---------- Row 17 ----------
$ws.Range("C17").Value = 61

# ---------- Row 18 ----------
$ws.Range("C18").Value = 75
Set-TextValue "D18" "31"
$ws.Range("E18").Value = "Double pole MCB(With B/C curve tripping Characteristics)"

# ---------- Row 19 ----------
Set-TextValue "A19" "Each"
$ws.Range("C19").Value = 85
Set-TextValue "D19" "35"
$ws.Range("E19").Value = "8 Way (8+2)"
$ws.Range("F19").Value = 2184
Set-TextValue "G19" "185640.00"

# ---------- Row 20 ----------
Set-TextValue "A20" "%"
$ws.Range("C20").Value = 51
Set-TextValue "D20" "37"
$ws.Range("E20").Value = "Add Tender Premium "
$ws.Range("F20").Value = 0
Set-TextValue "G20" "0.00"

# ---------- Delete row 21 (the old "Grand Total" item-row) ----------
# This shifts old rows 22-25 up to become new rows 21-24, matching the target layout.
$ws.Rows(21).Delete()

# ---------- Update the recomputed totals on the now-shifted summary rows ----------
# New row 22: "Grand Total Rs." row -> updated total
Set-TextValue "G22" "570563.00"
Set-TextValue "H22" "570563.00"

# New row 24: "NET PAYABLE AMOUNT Rs." row -> updated total
Set-TextValue "G24" "570563.00"
Set-TextValue "H24" "570563.00"
